$d = $word.ActiveDocument

function Get-ParagraphEnd($pos) {
    # Robustly resolve the content-end (pilcrow-exclusive) position of the
    # paragraph that contains document position $pos, using only the
    # document-level Paragraphs collection (Range.Paragraphs on a
    # sub-range / zero-width range is unreliable in this engine).
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $s = $p.Range.Start
        $e = $p.Range.End
        if (($pos -ge $s) -and ($pos -le $e)) {
            return $e
        }
    }
    throw "Could not resolve paragraph for position $pos"
}

# -----------------------------------------------------------------
# Edit 1: Add "Hossain" to the team-members list ("...Aaron" -> "...Aaron Hossain")
# Done as a brand-new trailing run by inserting at the exact end of the
# matched text (a reliable run-boundary in this engine, since "Aaron" is
# the paragraph's last word), mirroring how the original "Aaron" run
# already sits at the paragraph's end.
# -----------------------------------------------------------------
$namesRng = $d.Content
$foundNames = $namesRng.Find.Execute("Aaron", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundNames) {
    throw "Could not find 'Aaron' in the team members paragraph"
}
$namesInsertPos = $namesRng.End
$insertRng = $d.Range($namesInsertPos, $namesInsertPos)
$insertRng.InsertAfter(" Hossain")

# -----------------------------------------------------------------
# Edit 2: Touch up the "Jeremiah's jungle" description paragraph.
#   - insert "game " before "mechanics and movement"
#   - replace "for" with "to reclaim" before "Jeremiah's crown"
#   - insert "(TBD) " before "dropping coconuts"
#   - wrap "similar to" with gramStart/gramEnd proofErr markers
# The whole trailing run (from "The first level..." through the end of
# the paragraph) is rebuilt in one shot via Range.InsertXML so we get
# precise control over run boundaries and the w:proofErr elements.
# -----------------------------------------------------------------
$descRng = $d.Content
$foundDesc = $descRng.Find.Execute("The first level is a tutorial", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundDesc) {
    throw "Could not find the 'The first level is a tutorial' run"
}
$descRunStart = $descRng.Start
$descParaEnd = Get-ParagraphEnd $descRunStart

$targetRng = $d.Range($descRunStart, $descParaEnd)

$xmlSnippet = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">The first level is a tutorial, letting the player familiarize themselves with the </w:t></w:r><w:r><w:t xml:space="preserve">game </w:t></w:r><w:r><w:t xml:space="preserve">mechanics and movement, letting the player start on a stationary platform trying to navigate through moving and stationary platforms. After the first level, the player is rewarded with a dash ability they could use in the air and on the ground, the second level will allow the player to experiment with the new dash ability with smaller platforms, wider gaps and more obstacles to traverse through. The final level will be a boss fight </w:t></w:r><w:r><w:t>to reclaim</w:t></w:r><w:r><w:t xml:space="preserve"> Jeremiah’s crown, with a mysterious animal </w:t></w:r><w:r><w:t xml:space="preserve">(TBD) </w:t></w:r><w:r><w:t xml:space="preserve">dropping coconuts on Jeremiah, </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>similar to</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> Donkey Kong’s role in Donkey Kong. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$targetRng.InsertXML($xmlSnippet)
